$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing attendance status cells C2:C6 from "Presente" to "No vino"
$ws.Range("C2:C6").Value = "No vino"

# Add a new row for the 5th student (row 7)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "LeydiMoralesRubiano"
$ws.Range("C7").Value = "No vino"

# Copy the style from A6 (index/number column) to A7 so it matches formatting
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
